$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-05-27 Tuesday"; new="2025-05-28 Wednesday"},
    @{old="68×45=3060"; new="19×32=608"},
    @{old="60×72=4320"; new="50×44=2200"},
    @{old="77×83=6391"; new="62×54=3348"},
    @{old="97×96=9312"; new="36×25=900"},
    @{old="86×79=6794"; new="20×60=1200"},
    @{old="41×87=3567"; new="64×48=3072"},
    @{old="76×38=2888"; new="94×35=3290"},
    @{old="27×68=1836"; new="68×83=5644"},
    @{old="11×23=253"; new="70×31=2170"},
    @{old="36×37=1332"; new="90×29=2610"},
    @{old="56×68=3808"; new="20×14=280"},
    @{old="68×43=2924"; new="23×51=1173"},
    @{old="41×72=2952"; new="71×36=2556"},
    @{old="30×81=2430"; new="26×94=2444"},
    @{old="20×43=860"; new="22×94=2068"},
    @{old="39×12=468"; new="69×67=4623"},
    @{old="64×34=2176"; new="84×31=2604"},
    @{old="61×33=2013"; new="53×66=3498"},
    @{old="58×79=4582"; new="60×38=2280"},
    @{old="47×76=3572"; new="16×85=1360"},
    @{old="20×25=500"; new="85×81=6885"},
    @{old="70×76=5320"; new="99×86=8514"},
    @{old="68×16=1088"; new="36×63=2268"},
    @{old="16×29=464"; new="18×88=1584"},
    @{old="16×41=656"; new="11×54=594"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
